$d = $word.ActiveDocument

# The abstract body is the second paragraph. Replace its whole Range.Text
# (rather than doing piecemeal Find/Replace) so the several runs that make
# up the paragraph - including the isolated "á" run with its own rPr -
# collapse back down into the single plain run the edit settled on.
$para = $d.Paragraphs.Item(2)
$start = $para.Range.Start
$end = $para.Range.End
$r = $d.Range($start, $end)

$r.Text = "For a long time have models of the Moho discontinuity have been created from a variety of different methods including gravitationally and seismologically. However, very few of these models developed have uncertainty estimates, this is especially the case where there is limited seismic data available. In areas such as South America and Africa due to the economic and environmental challenges, over vast regions of the continent, there are little to no seismic point estimates. For these areas to have relevant Moho models either gravitational data needs to be used or the seismic data has to be interpolated but there is no way to tell how accurate these methods are to attaining a true Moho depth model. To get around this problem a method of cross-validation specifically repeated random sub-sample validation will be used to quantify errors on a gravitationally derived model of South America with the help of seismic point estimates. The results from this cross-validation will tell how good gravitational models are in regions where there is no seismic data to compare it to. Additionally, for regions where the model significantly underestimates the Moho depth in comparison to the seismic data, there is likely an unmodelled mass present. The Paraná Basin is thought to have large igneous intrusions resulting in a shallower Moho than expected, in an attempt to decrease the errors on the model this intrusion will be modelled. The cross-validation gives an indication that the specific model used has a reasonably small misfit from the point estimates but that the size of the error may vary geographically across South America as the error values attained are an average for the whole model."
